# IST price update 2025-12-20 18:07
#
# A new price-scrape column is inserted immediately to the right of the
# "SKU Name" column (i.e. the new, second column). The previous B/C/D
# price columns shift one place to the right (B->C, C->D, D->E).
# Row 1 holds the timestamp headers; the newest timestamp goes into the
# freshly inserted B1. Row 13's new scrape failed, so its new B cell is
# left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 26

# Shift existing D/C/B price values one column to the right (process the
# rightmost column first so values aren't clobbered before they are read).
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 2).Value2
}

# Populate the newly freed-up column B with the latest scrape.
$ws.Cells.Item(1, 2).Value = "2025-12-20 23:33"
$ws.Cells.Item(13, 2).Value = ""

# New column E should look like the others: same width as column D, and
# (for the header row) the same bold/bordered style as column D's header.
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
